$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates: Volume/Number and report week date range
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# Weekly crime statistics data updates (rows 14-30, columns C-N)
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 8
$ws.Range("E14").Value = 12.5
$ws.Range("F14").Value = 30
$ws.Range("G14").Value = 27
$ws.Range("H14").Value = 11.111111111111
$ws.Range("I14").Value = 391
$ws.Range("J14").Value = 440
$ws.Range("K14").Value = -11.136363636363
$ws.Range("L14").Value = -10.114942528735
$ws.Range("M14").Value = -20.366598778004
$ws.Range("N14").Value = -77.695379349686
$ws.Range("C15").Value = 16
$ws.Range("D15").Value = 29
$ws.Range("E15").Value = -44.827586206896
$ws.Range("F15").Value = 103
$ws.Range("G15").Value = 116
$ws.Range("H15").Value = -11.206896551724
$ws.Range("I15").Value = 1487
$ws.Range("J15").Value = 1365
$ws.Range("K15").Value = 8.937728937728
$ws.Range("L15").Value = 11.804511278195
$ws.Range("M15").Value = 18.675179569034
$ws.Range("N15").Value = -50.067159167226
$ws.Range("C16").Value = 317
$ws.Range("D16").Value = 293
$ws.Range("E16").Value = 8.191126279863
$ws.Range("F16").Value = 1417
$ws.Range("G16").Value = 1346
$ws.Range("H16").Value = 5.274888558692
$ws.Range("I16").Value = 15972
$ws.Range("J16").Value = 12326
$ws.Range("K16").Value = 29.579750121694
$ws.Range("L16").Value = 34.727962884858
$ws.Range("M16").Value = -9.563444878545
$ws.Range("N16").Value = -79.451684699404
$ws.Range("C17").Value = 441
$ws.Range("D17").Value = 420
$ws.Range("E17").Value = 5
$ws.Range("F17").Value = 1899
$ws.Range("G17").Value = 1870
$ws.Range("H17").Value = 1.550802139037
$ws.Range("I17").Value = 23747
$ws.Range("J17").Value = 20983
$ws.Range("K17").Value = 13.172568269551
$ws.Range("L17").Value = 24.938180670279
$ws.Range("M17").Value = 51.525012761613
$ws.Range("N17").Value = -37.532552939629
$ws.Range("C18").Value = 273
$ws.Range("D18").Value = 316
$ws.Range("E18").Value = -13.607594936708
$ws.Range("F18").Value = 1123
$ws.Range("G18").Value = 1189
$ws.Range("H18").Value = -5.550883095037
$ws.Range("I18").Value = 14281
$ws.Range("J18").Value = 11351
$ws.Range("K18").Value = 25.812703726543
$ws.Range("L18").Value = 1.607968694414
$ws.Range("M18").Value = -16.147026011391
$ws.Range("N18").Value = -84.375957288521
$ws.Range("C19").Value = 853
$ws.Range("D19").Value = 1004
$ws.Range("E19").Value = -15.039840637450
$ws.Range("F19").Value = 3900
$ws.Range("G19").Value = 3969
$ws.Range("H19").Value = -1.738473167044
$ws.Range("I19").Value = 46821
$ws.Range("J19").Value = 34923
$ws.Range("K19").Value = 34.069238037969
$ws.Range("L19").Value = 43.834480216269
$ws.Range("M19").Value = 36.159013580713
$ws.Range("N19").Value = -39.916845252608
$ws.Range("C20").Value = 255
$ws.Range("D20").Value = 262
$ws.Range("E20").Value = -2.671755725190
$ws.Range("F20").Value = 1047
$ws.Range("G20").Value = 975
$ws.Range("H20").Value = 7.384615384615
$ws.Range("I20").Value = 12297
$ws.Range("J20").Value = 9339
$ws.Range("K20").Value = 31.673626726630
$ws.Range("L20").Value = 50.054911531421
$ws.Range("M20").Value = 30.777411464426
$ws.Range("N20").Value = -87.931694391285
$ws.Range("C21").Value = 2164
$ws.Range("D21").Value = 2332
$ws.Range("E21").Value = -7.204116638078
$ws.Range("F21").Value = 9519
$ws.Range("G21").Value = 9492
$ws.Range("H21").Value = 0.284450063211
$ws.Range("I21").Value = 114996
$ws.Range("J21").Value = 90727
$ws.Range("K21").Value = 26.749479206851
$ws.Range("L21").Value = 31.53072779055
$ws.Range("M21").Value = 19.914909591440
$ws.Range("N21").Value = -70.641892668132
$ws.Range("C22").Value = 38
$ws.Range("D22").Value = 49
$ws.Range("E22").Value = -22.448979591836
$ws.Range("F22").Value = 185
$ws.Range("G22").Value = 213
$ws.Range("H22").Value = -13.145539906103
$ws.Range("I22").Value = 2096
$ws.Range("J22").Value = 1570
$ws.Range("K22").Value = 33.503184713375
$ws.Range("L22").Value = 28.826060233558
$ws.Range("M22").Value = 8.376421923474
$ws.Range("C23").Value = 110
$ws.Range("D23").Value = 103
$ws.Range("E23").Value = 6.796116504854
$ws.Range("G23").Value = 448
$ws.Range("H23").Value = 2.678571428571
$ws.Range("I23").Value = 5422
$ws.Range("J23").Value = 5009
$ws.Range("K23").Value = 8.245158714314
$ws.Range("L23").Value = 15.978609625668
$ws.Range("M23").Value = 40.612033195020
$ws.Range("C24").Value = 1845
$ws.Range("D24").Value = 1899
$ws.Range("E24").Value = -2.843601895734
$ws.Range("F24").Value = 8630
$ws.Range("G24").Value = 7933
$ws.Range("H24").Value = 8.786083448884
$ws.Range("I24").Value = 104726
$ws.Range("J24").Value = 77184
$ws.Range("K24").Value = 35.683561359867
$ws.Range("L24").Value = 41.759164004548
$ws.Range("M24").Value = 40.675666599503
$ws.Range("C25").Value = 704
$ws.Range("D25").Value = 727
$ws.Range("E25").Value = -3.163686382393
$ws.Range("F25").Value = 3113
$ws.Range("G25").Value = 3113
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 37696
$ws.Range("J25").Value = 32905
$ws.Range("K25").Value = 14.560097249658
$ws.Range("L25").Value = 23.865540695955
$ws.Range("M25").Value = -10.409734765662
$ws.Range("C26").Value = 30
$ws.Range("E26").Value = -41.176470588235
$ws.Range("F26").Value = 165
$ws.Range("G26").Value = 205
$ws.Range("H26").Value = -19.512195121951
$ws.Range("I26").Value = 2380
$ws.Range("J26").Value = 2245
$ws.Range("K26").Value = 6.013363028953
$ws.Range("L26").Value = 14.754098360655
$ws.Range("C27").Value = 74
$ws.Range("D27").Value = 86
$ws.Range("E27").Value = -13.953488372093
$ws.Range("F27").Value = 388
$ws.Range("G27").Value = 398
$ws.Range("H27").Value = -2.512562814070
$ws.Range("I27").Value = 4747
$ws.Range("J27").Value = 4471
$ws.Range("K27").Value = 6.173115634086
$ws.Range("L27").Value = 36.017191977077
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 33
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 93
$ws.Range("G28").Value = 145
$ws.Range("H28").Value = -35.862068965517
$ws.Range("I28").Value = 1461
$ws.Range("J28").Value = 1721
$ws.Range("K28").Value = -15.107495642068
$ws.Range("L28").Value = -15.255220417633
$ws.Range("M28").Value = -10.860280658938
$ws.Range("N28").Value = -73.054223533751
$ws.Range("C29").Value = 21
$ws.Range("D29").Value = 32
$ws.Range("E29").Value = -34.375
$ws.Range("F29").Value = 80
$ws.Range("G29").Value = 121
$ws.Range("H29").Value = -33.884297520661
$ws.Range("I29").Value = 1205
$ws.Range("J29").Value = 1435
$ws.Range("K29").Value = -16.027874564459
$ws.Range("L29").Value = -14.173789173789
$ws.Range("M29").Value = -10.806809770540
$ws.Range("N29").Value = -75.276979893311
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = 33.333333333333
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = 18.181818181818
$ws.Range("I30").Value = 570
$ws.Range("J30").Value = 484
$ws.Range("K30").Value = 17.768595041322
$ws.Range("L30").Value = 132.65306122449
